# Applies the "Add files via upload" commit:
#  - Sheet "DANH SÁCH NỢ" (1st sheet):
#      * H2: 25% -> 30% interest rate (F2/I2 formulas recalc automatically)
#      * Rows 5-7: mark as settled/highlighted (copy the already-used
#        "fillId 3" formatting from row 4) and clear the "Kết thúc" (K) date
#      * Row 25: fill in a brand-new debt entry (was a blank template row)
#      * View: scroll/zoom/selection tweaks
#  - Sheet "THONG KE NAP " (2nd sheet):
#      * Row 140: fill in the matching "nạp quân huy" transaction entry
#      * View: selection tweak

$wb = $excel.ActiveWorkbook

$wsDebt = $wb.Worksheets.Item(1)   # "DANH SÁCH NỢ"
$wsNap  = $wb.Worksheets.Item(2)   # "THONG KE NAP "

# ---------------------------------------------------------------------
# 1) DANH SÁCH NỢ - row 2: interest rate 25% -> 30%
# ---------------------------------------------------------------------
$wsDebt.Range("H2").Value = 0.3

# ---------------------------------------------------------------------
# 2) DANH SÁCH NỢ - rows 5..7: these loans are now fully repaid, so they
#    get the same highlighted formatting already used on row 4 (style
#    set with fillId 3) and their "Kết thúc" (end date, column K) is
#    cleared out.
# ---------------------------------------------------------------------
$wsDebt.Range("B4:M4").Copy()
$wsDebt.Range("B5:M7").PasteSpecial(-4122)   # xlPasteFormats
$wsDebt.Range("K5:K7").ClearContents()

# ---------------------------------------------------------------------
# 3) DANH SÁCH NỢ - row 25: new debt entry
#    Nguyễn Huỳnh Anh Thư / Nạp quân huy / 80,000 / start 46026 / end 46032
# ---------------------------------------------------------------------
$wsDebt.Range("B25").Value = "Nguyễn Huỳnh Anh Thư"
$wsDebt.Range("C25").Value = "Nạp quân huy"
$wsDebt.Range("D25").Value = 80000
$wsDebt.Range("E25").Value = 0
$wsDebt.Range("F25").Formula = "=(D25+I25)-E25"
$wsDebt.Range("G25").Value = 0
$wsDebt.Range("H25").Value = 0
$wsDebt.Range("I25").Formula = "=D25*H25"
$wsDebt.Range("J25").Value = 46026
$wsDebt.Range("K25").Value = 46032
$wsDebt.Range("M25").Value = "Chưa trả đủ"

# ---------------------------------------------------------------------
# 4) THONG KE NAP - row 140: matching "nạp quân huy" entry
# ---------------------------------------------------------------------
$wsNap.Range("A140").Value = 46026
$wsNap.Range("B140").Value = "Nguyễn Huỳnh Anh Thư"
$wsNap.Range("C140").Value = 80000
$wsNap.Range("D140").Value = "Nạp quân huy"

# ---------------------------------------------------------------------
# 5) View state: scroll/zoom/selection
# ---------------------------------------------------------------------
$wsDebt.Range("A16").Select()
$wsDebt.Range("H31").Select()
$excel.ActiveWindow.Zoom = 85

$wsNap.Range("D141").Select()

# Restore the originally active sheet/tab before finishing.
$wsDebt.Activate()
